$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2104.2632
$ws.Range("J32").Value = 1338.8462
$ws.Range("L32").Value = 1338.8462
$ws.Range("N32").Value = -1990.8462

# Hunk 1: ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 20002570
$ws.Range("I116").Value = 100001000
$ws.Range("J116").Value = 2962.375
$ws.Range("K116").Value = 100001000
$ws.Range("L116").Value = 2962.375
$ws.Range("M116").Value = -99997558
$ws.Range("N116").Value = -9846.375

# Hunk 2: ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 73779
$ws.Range("J133").Value = 73779
$ws.Range("L133").Value = 73779
$ws.Range("N133").Value = -83899

# Hunk 3: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 12212.671
$ws.Range("I135").Value = 14872.203
$ws.Range("J135").Value = 743.4375
$ws.Range("K135").Value = 133849.827
$ws.Range("L135").Value = 6690.9375
$ws.Range("M135").Value = -131314.827
$ws.Range("N135").Value = -11760.9375

# Hunk 4: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1670.75
$ws.Range("I137").Value = 1644.3334
$ws.Range("K137").Value = 4933.0002
$ws.Range("M137").Value = -2383.0002

# Hunk 5: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1932.6451
$ws.Range("J138").Value = 2705.6086
$ws.Range("L138").Value = 8116.825800000001
$ws.Range("N138").Value = -18396.8258

# Hunk 6: ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2387
$ws.Range("I2").Value = 2211
$ws.Range("J2").Value = 2475
$ws.Range("K2").Value = 2211
$ws.Range("L2").Value = 2475
$ws.Range("M2").Value = -2098
$ws.Range("N2").Value = -2701

# Hunk 7: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14697.813
$ws.Range("I32").Value = 15600.833
$ws.Range("J32").Value = 10053.714
$ws.Range("K32").Value = 15600.833
$ws.Range("L32").Value = 10053.714
$ws.Range("M32").Value = -15313.833
$ws.Range("N32").Value = -10627.714

# Hunk 8: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2174.6
$ws.Range("I74").Value = 2174.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2174.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1300.6
$ws.Range("N74").ClearContents()

# Hunk 9: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2174.6
$ws.Range("I77").Value = 2174.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 10873
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -6505
$ws.Range("N77").ClearContents()

# Hunk 10: ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3550
$ws.Range("I88").Value = 2270
$ws.Range("K88").Value = 2270
$ws.Range("M88").Value = -1864

# Hunk 11: ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3550
$ws.Range("I91").Value = 2270
$ws.Range("K91").Value = 2270
$ws.Range("M91").Value = -866

# Hunk 12: ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 275000
$ws.Range("J103").Value = 275000
$ws.Range("L103").Value = 275000
$ws.Range("N103").Value = -277344

# Hunk 13: ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2387
$ws.Range("I116").Value = 2211
$ws.Range("J116").Value = 2475
$ws.Range("K116").Value = 2211
$ws.Range("L116").Value = 2475
$ws.Range("M116").Value = 83
$ws.Range("N116").Value = -7063

# Hunk 14: BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2387
$ws.Range("I3").Value = 2211
$ws.Range("J3").Value = 2475
$ws.Range("K3").Value = 2211
$ws.Range("L3").Value = 2475
$ws.Range("M3").Value = -2097
$ws.Range("N3").Value = -2703

# Hunk 15: BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 145902.86
$ws.Range("J86").Value = 203384
$ws.Range("L86").Value = 203384
$ws.Range("N86").Value = -205630

# Hunk 16: BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 145902.86
$ws.Range("J89").Value = 203384
$ws.Range("L89").Value = 1016920
$ws.Range("N89").Value = -1028152

# Hunk 17: BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 48000
$ws.Range("J126").Value = 48000
$ws.Range("L126").Value = 48000
$ws.Range("N126").Value = -57880

# Hunk 18: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2198.5588
$ws.Range("I134").Value = 1628.4348
$ws.Range("J134").Value = 3390.6365
$ws.Range("K134").Value = 4885.3044
$ws.Range("L134").Value = 10171.9095
$ws.Range("M134").Value = -2350.3044
$ws.Range("N134").Value = -15241.9095

# Hunk 19: CRP row 28
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 167821.5
$ws.Range("J28").Value = 167821.5
$ws.Range("L28").Value = 167821.5
$ws.Range("N28").Value = -168311.5

# Hunk 20: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3130.1177
$ws.Range("I31").Value = 1839
$ws.Range("K31").Value = 1839
$ws.Range("M31").Value = -1544

# Hunk 21: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3130.1177
$ws.Range("I34").Value = 1839
$ws.Range("K34").Value = 1839
$ws.Range("M34").Value = -1637

# Hunk 22: CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 102621
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# Hunk 23: CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 102621
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# Hunk 24: CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 19750
$ws.Range("J95").Value = 19750
$ws.Range("L95").Value = 19750
$ws.Range("N95").Value = -25242

# Hunk 25: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3310.6875
$ws.Range("I132").Value = 2437.3
$ws.Range("J132").Value = 4766.3335
$ws.Range("K132").Value = 7311.900000000001
$ws.Range("L132").Value = 14299.0005
$ws.Range("M132").Value = -4781.900000000001
$ws.Range("N132").Value = -19359.0005

# Hunk 26: CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 905.2308
$ws.Range("I2").Value = 1776.3334
$ws.Range("J2").Value = 158.57143
$ws.Range("K2").Value = 10658.0004
$ws.Range("L2").Value = 951.42858
$ws.Range("M2").Value = -10545.0004
$ws.Range("N2").Value = -1177.42858

# Hunk 27: CUL row 15
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 400
$ws.Range("J15").Value = 400
$ws.Range("L15").Value = 1200
$ws.Range("N15").Value = -1480

# Hunk 28: CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Hunk 29: CUL row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 4274.4546
$ws.Range("I36").Value = 3003.2222
$ws.Range("J36").Value = 9995
$ws.Range("K36").Value = 9009.6666
$ws.Range("L36").Value = 29985
$ws.Range("M36").Value = -8840.6666
$ws.Range("N36").Value = -30323

# Hunk 30: CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 628.53656
$ws.Range("I113").Value = 595.21875
$ws.Range("K113").Value = 1785.65625
$ws.Range("M113").Value = 384.34375

# Hunk 31: GSM row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 25000
$ws.Range("J98").Value = 25000
$ws.Range("L98").Value = 25000
$ws.Range("N98").Value = -30990

# Hunk 32: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3162.35
$ws.Range("I132").Value = 2128.7144
$ws.Range("J132").Value = 3718.923
$ws.Range("K132").Value = 6386.1432
$ws.Range("L132").Value = 11156.769
$ws.Range("M132").Value = -3856.1432
$ws.Range("N132").Value = -16216.769

# Hunk 33: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4337.5
$ws.Range("J7").Value = 3400
$ws.Range("L7").Value = 3400
$ws.Range("N7").Value = -3624

# Hunk 34: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4337.5
$ws.Range("J126").Value = 3400
$ws.Range("L126").Value = 10200
$ws.Range("N126").Value = -15140

# Hunk 35: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5692.4287
$ws.Range("I132").Value = 5069.4
$ws.Range("J132").Value = 7250
$ws.Range("K132").Value = 15208.2
$ws.Range("L132").Value = 21750
$ws.Range("M132").Value = -12678.2
$ws.Range("N132").Value = -26810

# Hunk 36: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 25252102
$ws.Range("I136").Value = 32259390
$ws.Range("J136").Value = 1115888.4
$ws.Range("K136").Value = 96778170
$ws.Range("L136").Value = 3347665.2
$ws.Range("M136").Value = -96775620
$ws.Range("N136").Value = -3352765.2

# Hunk 37: WVR row 97
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 40571.5
$ws.Range("J97").Value = 40571.5
$ws.Range("L97").Value = 40571.5
$ws.Range("N97").Value = -42553.5

# Hunk 38: WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 591.8570999999999
$ws.Range("I107").Value = 386.66666
$ws.Range("J107").Value = 745.75
$ws.Range("K107").Value = 1159.99998
$ws.Range("L107").Value = 2237.25
$ws.Range("M107").Value = 760.0000199999999
$ws.Range("N107").Value = -6077.25

# Hunk 39: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2268.8333
$ws.Range("I132").Value = 1161.75
$ws.Range("K132").Value = 3485.25
$ws.Range("M132").Value = -955.25

# Hunk 40: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1570.909
$ws.Range("I136").Value = 1716.3334
$ws.Range("J136").Value = 1259.2858
$ws.Range("K136").Value = 5149.0002
$ws.Range("L136").Value = 3777.8574
$ws.Range("M136").Value = -2599.0002
$ws.Range("N136").Value = -8877.857400000001
